$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 14th place as a new row at the bottom of the table
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Yerevan, Artsakh Street"
$ws.Range("C15").Value = "Apartment"
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 40.1447873580964
$ws.Range("F15").Value = 44.5100795865558
$ws.Range("G15").Value = "Erebuni"

# Update the active selection as recorded after the edit
$ws.Range("B19").Select()
